$d = $word.ActiveDocument

# Replace the identifier code (appears in two places)
$d.Content.Find.Execute("6.11.2.11.3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2.28.2.1.8", 2)

# Replace the document path reference (appears in two places)
$d.Content.Find.Execute("KUR.0130.00UNZ.SBA.TS.PA0046", $true, $false, $false, $false, $false,
                         $true, 1, $false, "KUR.0130.00USY.0.TZ.PA0025", 2)
